# reference-errors.xlsx fixture update:
#  - "Main root" sheet header A1 renamed from "Id" to "Identifier" (exercises
#    the new inexact / case-insensitive attribute-name matching added in this
#    commit).
#  - Active selection on the "Nodes" sheet moved to B3 (incidental editor
#    state captured when the fixture was last saved).

$wb = $excel.ActiveWorkbook

$mainRoot = $wb.Worksheets.Item("Main root")
$mainRoot.Range("A1").Value = "Identifier"

$nodes = $wb.Worksheets.Item("Nodes")
$nodes.Range("B3").Select()
